function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cur = $shapes.Item($i)
        if ($cur.Id -eq $id) {
            return $cur
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s8 = $p.Slides.Item(8)

# Shape id=14 "Rectangle 13" - HelpWindow -> WelcomeWindow, widened
$rect13 = Get-ShapeById $s8.Shapes 14
$rect13.TextFrame.TextRange.Text = "WelcomeWindow"
$rect13.Width = 1185420

# Shape id=21 "Elbow Connector 20" - height nudge (auto re-route)
$conn20 = Get-ShapeById $s8.Shapes 21
$conn20.Height = 420378

# Shape id=28 "Elbow Connector 27" - re-routed position/size
$conn27 = Get-ShapeById $s8.Shapes 28
$conn27.Left = 3972601
$conn27.Top = 2809200
$conn27.Height = 1752021

# Shape id=31 "Rectangle 30" - remove spellcheck err flag is cosmetic; text unchanged ("Ui")

# Shape id=38 "Rectangle 37" - nudge down
$rect37 = Get-ShapeById $s8.Shapes 38
$rect37.Top = 2805240

# Shape id=40 "Rectangle 39" - nudge up
$rect39 = Get-ShapeById $s8.Shapes 40
$rect39.Top = 4411938
